$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row from column A (header in row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Header row: new columns I ("I0") and J ("IF") ---
# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they pick up the same style (bold/centered/bordered)
# instead of Excel creating a brand-new style entry.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows: I = constant 1 ("I0"), J = copy of H ("IP" -> "IF") ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ip = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ip
}
